# GQD_PHIEU_XAC_MINH_KHIEU_NAI_1.docx
#
# The document previously had two sentences split across three
# consecutive <w:r> runs each (because earlier edits had been done
# piecemeal). This change re-joins each sentence back into a single
# run with the same (unchanged) run formatting, without altering the
# wording at all.
#
# We do this with Find/Replace: searching for (and replacing with the
# same) the full sentence text causes Word to rewrite the matched
# range as a single run, collapsing the run boundaries that used to
# split the sentence.

$d = $word.ActiveDocument

# --- Paragraph: "... báo cáo kết quả xác minh nội dung khiếu nại nêu
#     trên với ${tenLanhDao} trước ngày... tháng ... năm ..."
# Previously split as:
#   "...nêu trên vớ" | "i ${tenLanhDao}" | " trước ngày... tháng ... năm ..."
$text1 = "`${tenCoQuanDuocGiaoNhiemVuXM} báo cáo kết quả xác minh nội dung khiếu nại nêu trên với `${tenLanhDao} trước ngày... tháng ... năm ..."
$found1 = $d.Content.Find.Execute($text1, $true, $false, $false, $false, $false, $true, 1, $false, $text1, 2)
if (-not $found1) {
    Write-Host "WARNING: paragraph 1 text not found"
}

# --- Paragraph: "Người đứng đầu ${tenCoQuanDuocGiaoNhiemVuXM},
#     ....…………………...(3) chịu trách nhiệm thi hành Quyết định này."
# Previously split as:
#   "Người đứng đầu ${tenCoQuanDuocGiaoNhiemVuXM}, ....…………………..." | "(3)" | " chịu trách nhiệm thi hành Quyết định này."
$text2 = "Người đứng đầu `${tenCoQuanDuocGiaoNhiemVuXM}, ....…………………...(3) chịu trách nhiệm thi hành Quyết định này."
$found2 = $d.Content.Find.Execute($text2, $true, $false, $false, $false, $false, $true, 1, $false, $text2, 2)
if (-not $found2) {
    Write-Host "WARNING: paragraph 2 text not found"
}

Write-Host "Paragraph 1 merged: $found1"
Write-Host "Paragraph 2 merged: $found2"
